$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New story card rows appended to the end of the table (rows 21-23)

# Row 21
$ws.Cells.Item(21, 1).Value = 14
$ws.Cells.Item(21, 2).Value = "jungfräulich"
$ws.Cells.Item(21, 3).Value = "mittel"
$ws.Cells.Item(21, 4).Value = "Öffentlichkeitsarbeit"
$ws.Cells.Item(21, 5).Value = "2011-10-07"
$ws.Cells.Item(21, 6).Value = "Salzer"
$ws.Cells.Item(21, 7).Value = "Wiederschein"
$ws.Cells.Item(21, 8).Value = "Recherche"
$ws.Cells.Item(21, 9).Value = "Sammeln von infos zu Ameisensysteme"

# Row 22
$ws.Cells.Item(22, 1).Value = 15
$ws.Cells.Item(22, 2).Value = "jungfräulich"
$ws.Cells.Item(22, 3).Value = "mittel"
$ws.Cells.Item(22, 4).Value = "Öffentlichkeitsarbeit"
$ws.Cells.Item(22, 5).Value = "2011-10-07"
$ws.Cells.Item(22, 6).Value = "Salzer"
$ws.Cells.Item(22, 7).Value = "Wiederschein"
$ws.Cells.Item(22, 8).Value = "Recherche"
$ws.Cells.Item(22, 9).Value = "Sammeln von infos zu TSP"

# Row 23
$ws.Cells.Item(23, 1).Value = 16
$ws.Cells.Item(23, 2).Value = "jungfräulich"
$ws.Cells.Item(23, 3).Value = "hoch"
$ws.Cells.Item(23, 4).Value = "Öffentlichkeitsarbeit"
$ws.Cells.Item(23, 5).Value = "2011-10-07"
$ws.Cells.Item(23, 6).Value = "Salzer"
$ws.Cells.Item(23, 7).Value = "Wiederschein"
$ws.Cells.Item(23, 8).Value = "Ausarbeitung"
$ws.Cells.Item(23, 9).Value = "Komprimierung der gesammelten Daten"

# Reflect the final cursor position left by the edit session
$ws.Range("G26").Select()
